# Add a new logged problem ("Merge Two Sorted Linked Lists") as row 3
# on the "Linked List" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Linked List")

# Column A holds a date formatted as plain text (e.g. "08/09/2025" in
# row 2 above). Force the cell to Text format before assigning so Excel
# doesn't auto-convert the "MM/DD/YYYY" string into a real date serial,
# then drop back to the workbook's default (Normal) style so no stray
# cell formatting is introduced.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "08/11/2025"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "Merge Two Sorted Linked Lists"
$ws.Range("C3").Value = "Linked List"
$ws.Range("D3").Value = "Easy"
$ws.Range("E3").Value = "No (knew method but didnt know how to code it)"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = "Yes"
